# CRM-1792 Add RM detail in SF daily report
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells (row 1) and sample/placeholder cells (row 2) for the
# Regional Manager fields appended at columns AE/AF.
$ws.Range("AE1").Value = "Regional Manager Name"
$ws.Range("AF1").Value = "Regional Manager Contact Name"
$ws.Range("AE2").Value = "{vendor:sf_rm_name}"
$ws.Range("AF2").Value = "{vendor:sf_rm_phone}"

# Column width adjustments: AE/AF no longer share one width, each gets its
# own slightly wider size to fit the new labels.
$ws.Columns.Item(31).ColumnWidth = 22.77734375
$ws.Columns.Item(32).ColumnWidth = 27.21875

# Update the window scroll/selection state to match the saved view.
$ws.Application.ActiveWindow.ScrollColumn = 25
$ws.Range("AD11").Select()
